# Apply the "F column" (想去人数) updates across sheets, and remove row 2
# from the 本地生活 sheet (the single extra event row).

$wb = $excel.ActiveWorkbook

# --- Sheet 展览 (Sheet1): update F column values ---
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    3  = 39
    5  = 154
    6  = 1005
    8  = 7790
    11 = 6769
    12 = 148
    14 = 4807
    15 = 5170
    16 = 1058
    17 = 298
    18 = 302
    26 = 8813
    28 = 1562
    29 = 36
    31 = 796
    32 = 67
    33 = 139
    34 = 999
    37 = 1108
    39 = 4598
    40 = 21
    41 = 369
    43 = 47
    44 = 135
    46 = 21
    47 = 1206
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# --- Sheet 演出 (Sheet2): update F column values ---
$ws2 = $wb.Worksheets.Item("演出")
$updates2 = @{
    4  = 5
    17 = 881
}
foreach ($row in $updates2.Keys) {
    $ws2.Range("F$row").Value = $updates2[$row]
}

# --- Sheet 本地生活 (Sheet3): delete row 2 (the only data row) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Rows.Item(2).Delete()

# --- Sheet 全部类型 (Sheet4): update F column values ---
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    4  = 39
    6  = 154
    8  = 1005
    10 = 7790
    13 = 6769
    14 = 148
    16 = 4807
    17 = 5170
    18 = 1058
    19 = 298
    20 = 302
    27 = 8813
    29 = 1562
    30 = 36
    32 = 796
    33 = 67
    34 = 139
    35 = 999
    37 = 1108
    39 = 4598
    40 = 21
    41 = 369
    43 = 47
    44 = 135
    46 = 21
    47 = 1206
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
